$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set the Price column (D) cells to Text format first so that numeric-looking
# strings (e.g. "19.70", "1.00") are preserved exactly as text, matching the
# original inline-string cells, instead of being coerced into numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "70.673.82"
$ws.Range("E2").Value = "  +1.23%  "

$ws.Range("D3").Value = "3.581.98"

$ws.Range("E4").Value = "  +0.01%  "

$ws.Range("D5").Value = "589.13"
$ws.Range("E5").Value = "  +2.51%  "

$ws.Range("D6").Value = "187.29"
$ws.Range("E6").Value = "  +1.11%  "

$ws.Range("D7").Value = "3.570.51"
$ws.Range("E7").Value = "  +0.46%  "

$ws.Range("E8").Value = "  +0.81%  "

$ws.Range("E9").Value = "  +0.02%  "

$ws.Range("D10").Value = "0.201"
$ws.Range("E10").Value = "  +10.17%  "

$ws.Range("D11").Value = "0.651"
$ws.Range("E11").Value = "  +0.83%  "

$ws.Range("D12").Value = "54.61"
$ws.Range("E12").Value = "  -0.12%  "

$ws.Range("D13").Value = "0.0000312"
$ws.Range("E13").Value = "  +4.03%  "

$ws.Range("D14").Value = "9.57"
$ws.Range("E14").Value = "  +0.92%  "

$ws.Range("D15").Value = "4.153.23"
$ws.Range("E15").Value = "  +0.58%  "

$ws.Range("D16").Value = "19.70"
$ws.Range("E16").Value = "  +1.04%  "

$ws.Range("D17").Value = "70.667.26"
$ws.Range("E17").Value = "  +1.27%  "

$ws.Range("D18").Value = "3.585.49"
$ws.Range("E18").Value = "  +0.85%  "

$ws.Range("E19").Value = "  -0.12%  "

$ws.Range("B20").Value = "TRON"
$ws.Range("C20").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D20").Value = "0.120"
$ws.Range("E20").Value = "  +0.02%  "

$ws.Range("B21").Value = "BitcoinCash"
$ws.Range("C21").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D21").Value = "561.81"
$ws.Range("E21").Value = "  +14.43%  "

$ws.Range("E22").Value = "  -0.52%  "

$ws.Range("D23").Value = "17.91"
$ws.Range("E23").Value = "  -6.90%  "

$ws.Range("D24").Value = "4.69"
$ws.Range("E24").Value = "  +7.72%  "

$ws.Range("D25").Value = "4.92"
$ws.Range("E25").Value = "  +1.21%  "

$ws.Range("D26").Value = "95.85"
$ws.Range("E26").Value = "  +1.02%  "

$ws.Range("E27").Value = "  +1.60%  "

$ws.Range("D28").Value = "2.98"
$ws.Range("E28").Value = "  +1.66%  "

$ws.Range("D29").Value = "9.17"
$ws.Range("E29").Value = "  -0.74%  "

$ws.Range("D30").Value = "32.27"
$ws.Range("E30").Value = "  +2.50%  "

$ws.Range("D31").Value = "7.29"
$ws.Range("E31").Value = "  -2.08%  "

$ws.Range("D32").Value = "12.50"
$ws.Range("E32").Value = "  +4.44%  "

$ws.Range("D33").Value = "65.13"
$ws.Range("E33").Value = "  -1.86%  "

$ws.Range("E34").Value = "  +0.70%  "

$ws.Range("D35").Value = "563.21"
$ws.Range("E35").Value = "  -0.26%  "

$ws.Range("D36").Value = "3.30"
$ws.Range("E36").Value = "  +4.79%  "

$ws.Range("E37").Value = "  +5.57%  "

$ws.Range("D38").Value = "38.11"
$ws.Range("E38").Value = "  -1.26%  "

$ws.Range("E39").Value = "  -0.05%  "

$ws.Range("D40").Value = "0.0₃0776"
$ws.Range("E40").Value = "  -1.20%  "

$ws.Range("E41").Value = "  +1.34%  "

$ws.Range("D42").Value = "3.350.25"
$ws.Range("E42").Value = "  +4.39%  "

$ws.Range("D43").Value = "3.36"
$ws.Range("E43").Value = "  -3.94%  "

$ws.Range("E44").Value = "  -2.46%  "

$ws.Range("E45").Value = "  +4.00%  "

$ws.Range("D46").Value = "2.98"
$ws.Range("E46").Value = "  +0.22%  "

$ws.Range("E47").Value = "  +1.79%  "

$ws.Range("D48").Value = "9.38"
$ws.Range("E48").Value = "  -0.90%  "

$ws.Range("E49").Value = "  +1.38%  "

$ws.Range("D50").Value = "1.00"
$ws.Range("E50").Value = "  +0.09%  "

$ws.Range("E51").Value = "  +19.60%  "

# Restore the default (Normal) style on column D so no stray text-format
# styling is left behind on cells that did not have one originally.
$ws.Range("D2:D51").Style = "Normal"
